$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
  2  = @(0.3529568509324034, 0.01031932399247814)
  3  = @(0.3753492847338557, 0.006713211876517052)
  4  = @(0.3280025554987386, 0.01898826826776655)
  5  = @(0.3608313878960405, 0.009451377997013214)
  6  = @(0.3467161657400654, 0.006767853111292094)
  7  = @(0.3445554540192456, 0.005710428290169664)
  8  = @(0.1278987591963732, 0.007671841152495148)
  9  = @(0.1963535850000913, 0.005959803150005001)
  10 = @(0.1934124930164705, 0.005405594801197919)
  11 = @(0.4336849683193303, 0.005586959350726499)
  12 = @(0.1701176367835098, 0.006247883150787478)
  13 = @(0.1671487427099686, 0.006594963417218175)
  14 = @(0.1946555227291147, 0.004008863048839792)
  15 = @(0.3685857288298887, 0.007377739101505094)
  16 = @(0.357503836925234,  0.008867140287160763)
  17 = @(0.1528193468450992, 0.004443970138827808)
  18 = @(0.34739219112078,   0.007779176928094804)
  19 = @(0.3625814555492692, 0.009162588567545244)
  20 = @(0.1756878579169599, 0.003934342509300924)
  21 = @(0.3206532438156445, 0.01263468466705323)
  22 = @(0.1204666841463396, 0.003929398309853801)
  23 = @(0.2308820415282557, 0.00739496739651905)
  24 = @(0.3166143844802502, 0.01245546723261897)
  25 = @(0.1031270709062008, 0.002964278980857023)
  26 = @(0.1393149614132256, 0.004284884607912248)
  27 = @(0.3527221061438221, 0.01011122822442783)
  28 = @(0.2886191377714448, 0.0100507336370167)
}

foreach ($row in $values.Keys) {
  $pair = $values[$row]
  $ws.Cells.Item($row, 2).Value = $pair[0]
  $ws.Cells.Item($row, 3).Value = $pair[1]
}
